$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Risk"
$ws.Range("C1").Value = "Curve"
$ws.Range("D1").Value = "Type"

$ws.Range("B2").Value = 637.4

$ws.Range("F8").WrapText = $true
$ws.Range("D7").Select()
